$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "weight" column before column D ("grade"), shifting
# grade -> E and infected -> F.
$ws.Columns.Item(4).Insert()

$ws.Range("D1").Value = "weight"

# Correct a couple of jump_length (column B) values that were re-measured.
$ws.Range("B5").Value = 5.6
$ws.Range("B6").Value = 9.1
$ws.Range("B7").Value = 8.2

# Fill in the new weight values for each row.
$weights = @(2.1, 2.3, 2.8, 2.4, 1.2, 4.1, 3.2, 1.1, 2.1, 2.4, 2.1, 1.5, 3.7, 2.9)
for ($i = 0; $i -lt $weights.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $weights[$i]
}
